# Quarterly indexing esoteric bug-fix operation
#
# Column A held the first-of-quarter date used to index each forecast row.
# The indexing was off: it should point to the 15th of the *following*
# month instead of the 1st of the quarter-start month. Shift every date
# in A2:A73 forward by one month and move the day-of-month to the 15th.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 73

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range("A$r")
    $oldDate = $cell.Value()
    $newDate = $oldDate.AddMonths(1).AddDays(14)
    $cell.Value = $newDate
}
